# Append the new "B2-A2" group schedule rows (61-101) to the sheet, mirroring
# the existing "B2-A1" rows' alternating row styling (odd rows reuse the style
# of row 59, even rows reuse the style of row 60).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("Year 5","B2-A2","endocrinology","1","06/12/2025","09:00:00",360),
    @("Year 5","B2-A2","endocrinology","2","07/12/2025","09:00:00",360),
    @("Year 5","B2-A2","endocrinology","3","08/12/2025","09:00:00",360),
    @("Year 5","B2-A2","endocrinology","4","09/12/2025","09:00:00",360),
    @("Year 5","B2-A2","endocrinology","5","10/12/2025","09:00:00",360),
    @("Year 5","B2-A2","endocrinology","6","13/12/2025","09:00:00",360),
    @("Year 5","B2-A2","endocrinology","7","14/12/2025","09:00:00",360),
    @("Year 5","B2-A2","endocrinology","8","15/12/2025","09:00:00",360),
    @("Year 5","B2-A2","endocrinology","9","16/12/2025","09:00:00",360),
    @("Year 5","B2-A2","endocrinology","10","17/12/2025","09:00:00",360),
    @("Year 5","B2-A2","gastroenterology","1","20/12/2025","09:00:00",360),
    @("Year 5","B2-A2","gastroenterology","2","21/12/2025","09:00:00",360),
    @("Year 5","B2-A2","gastroenterology","3","22/12/2025","09:00:00",360),
    @("Year 5","B2-A2","gastroenterology","4","23/12/2025","09:00:00",360),
    @("Year 5","B2-A2","gastroenterology","5","24/12/2025","09:00:00",360),
    @("Year 5","B2-A2","gastroenterology","6","27/12/2025","09:00:00",360),
    @("Year 5","B2-A2","gastroenterology","7","28/12/2025","09:00:00",360),
    @("Year 5","B2-A2","gastroenterology","8","29/12/2025","09:00:00",360),
    @("Year 5","B2-A2","gastroenterology","9","30/12/2025","09:00:00",360),
    @("Year 5","B2-A2","gastroenterology","10","31/12/2025","09:00:00",360),
    @("Year 5","B2-A2","nephrology","1","03/01/2026","09:00:00",360),
    @("Year 5","B2-A2","nephrology","2","04/01/2026","09:00:00",360),
    @("Year 5","B2-A2","nephrology","3","05/01/2026","09:00:00",360),
    @("Year 5","B2-A2","nephrology","4","06/01/2026","09:00:00",360),
    @("Year 5","B2-A2","nephrology","5","07/01/2026","09:00:00",360),
    @("Year 5","B2-A2","neurology","1","17/01/2026","09:00:00",360),
    @("Year 5","B2-A2","neurology","2","18/01/2026","09:00:00",360),
    @("Year 5","B2-A2","neurology","3","19/01/2026","09:00:00",360),
    @("Year 5","B2-A2","neurology","4","20/01/2026","09:00:00",360),
    @("Year 5","B2-A2","neurology","5","21/01/2026","09:00:00",360),
    @("Year 5","B2-A2","neurology","6","07/02/2026","09:00:00",360),
    @("Year 5","B2-A2","neurology","7","08/02/2026","09:00:00",360),
    @("Year 5","B2-A2","neurology","8","09/02/2026","09:00:00",360),
    @("Year 5","B2-A2","neurology","9","10/02/2026","09:00:00",360),
    @("Year 5","B2-A2","physical medicine","1","04/02/2026","09:00:00",360),
    @("Year 5","B2-A2","physical medicine","2","11/02/2026","09:00:00",360),
    @("Year 5","B2-A2","rheumatology","1","10/01/2026","09:00:00",360),
    @("Year 5","B2-A2","rheumatology","2","11/01/2026","09:00:00",360),
    @("Year 5","B2-A2","rheumatology","3","12/01/2026","09:00:00",360),
    @("Year 5","B2-A2","rheumatology","4","13/01/2026","09:00:00",360),
    @("Year 5","B2-A2","rheumatology","5","14/01/2026","09:00:00",360)
)

$startRow = 61
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    # Rows alternate the same way the existing data does: odd destination
    # rows copy the "odd" style template (row 59), even rows copy the "even"
    # style template (row 60) - this keeps the shading/number-format pattern
    # (and reuses the same style ids) instead of minting new ones.
    if ($r % 2 -eq 1) {
        $templateRow = 59
    } else {
        $templateRow = 60
    }

    $ws.Range("A$templateRow`:G$templateRow").Copy()
    $ws.Range("A$r`:G$r").PasteSpecial(-4122)  # xlPasteFormats

    $row = $newRows[$i]

    # Columns A-F are stored as text in the source workbook (inline strings),
    # including the numeric-looking Session/Date/Time values - a leading
    # apostrophe forces text entry instead of Excel's automatic
    # number/date coercion, while keeping the pasted display format.
    $ws.Cells.Item($r,1).Value = $row[0]
    $ws.Cells.Item($r,2).Value = $row[1]
    $ws.Cells.Item($r,3).Value = $row[2]
    $ws.Cells.Item($r,4).Value = "'" + $row[3]
    $ws.Cells.Item($r,5).Value = "'" + $row[4]
    $ws.Cells.Item($r,6).Value = "'" + $row[5]
    $ws.Cells.Item($r,7).Value = $row[6]
}

$excel.CutCopyMode = 0
